# Add the newly-created login entry (Admin / a#a) as a new row below the
# existing data, then leave the selection where the user last clicked (K15),
# matching the saved workbook view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 1).Value = "Admin"
$ws.Cells.Item(5, 2).Value = "a#a"

$ws.Range("K15").Select()
